# Add three more reporting days (columns O, P, Q) to the "CodeBook" sheet,
# mirroring the existing column N (the most-recent day at the time) and then
# patching the handful of cells whose counts actually moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bring over formatting (number format / font / style index) for the
#        new columns by copying the format of column N, which already has
#        the right cellXfs (s="5" for the data rows, s="1" for the totals
#        row) applied. xlPasteFormats = -4122.
$ws.Range("N3:N25").Copy() | Out-Null
$ws.Range("O3:Q25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 2. Seed O/P/Q (rows 3-24) with the same per-place counts as column N -
#        almost every value is unchanged day-over-day.
$ws.Range("O3:O24").Value2 = $ws.Range("N3:N24").Value2
$ws.Range("P3:P24").Value2 = $ws.Range("N3:N24").Value2
$ws.Range("Q3:Q24").Value2 = $ws.Range("N3:N24").Value2

# --- 3. Patch the cells whose cumulative counts actually increased on the
#        new days.
$ws.Range("O3").Value2 = 108
$ws.Range("P3").Value2 = 109
$ws.Range("Q3").Value2 = 111

$ws.Range("Q6").Value2 = 29

$ws.Range("P11").Value2 = 44
$ws.Range("Q11").Value2 = 44

# --- 4. Extend the totals row with SUM formulas over the new columns.
$ws.Range("O25").Formula = "=SUM(O3:O24)"
$ws.Range("P25").Formula = "=SUM(P3:P24)"
$ws.Range("Q25").Formula = "=SUM(Q3:Q24)"

# --- 5. Leave the selection where the author left it: the new last row of
#        totals.
$ws.Range("O25:Q25").Select() | Out-Null
